$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row with Roll = 9 ("Reach" column), Zombie column is empty -> add "Bone Tongue"
$cellBoneTongue = $table.Cell(10, 5)
$cellBoneTongue.Range.Text = "Bone Tongue"
$cellBoneTongue.Range.Font.Name = "Pericles"
$cellBoneTongue.Range.Font.Size = 8
$cellBoneTongue.Range.Font.SizeBi = 8

# Row with Roll = 10 ("Slam" column), Zombie column is empty -> add "Detect Life"
$cellDetectLife = $table.Cell(11, 5)
$cellDetectLife.Range.Text = "Detect Life"
$cellDetectLife.Range.Font.Name = "Pericles"
$cellDetectLife.Range.Font.Size = 8
$cellDetectLife.Range.Font.SizeBi = 8
